$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A14 with refined timestamp value (precision correction)
$ws.Cells.Item(14, 1).Value = 45877.54186517361

# Append new row 15 with the latest sensor reading
$ws.Cells.Item(15, 1).Value = 45877.58355226297
$ws.Cells.Item(15, 2).Value = 2025
$ws.Cells.Item(15, 3).Value = 32
$ws.Cells.Item(15, 4).Value = 19.39
$ws.Cells.Item(15, 5).Value = 77.70999999999999
$ws.Cells.Item(15, 6).Value = 66.8
$ws.Cells.Item(15, 7).Value = 13.51
$ws.Cells.Item(15, 8).Value = "SE"
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = "14:00:18"

# Match the date-formatted number format used by the rest of column A
$ws.Cells.Item(15, 1).NumberFormat = $ws.Cells.Item(14, 1).NumberFormat
